$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the md (H) formula so it computes D<row> - D$5 instead of D$5 - D<row>,
# which makes the md values correctly negative (and flips lcl/ucl accordingly).
$ws.Range("H2").Formula = "=D2-D`$5"
$ws.Range("H3").Formula = "=D3-D`$5"
$ws.Range("H4").Formula = "=D4-D`$5"

# Update the last-saved selection to F11 (matches the authored state).
$ws.Range("F11").Select()
